$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.384.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.27%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.891.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.33%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'237.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.16%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4838"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.55%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2902"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.90%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.77%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.895.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.05%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'16.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.03%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07391"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.65%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.176"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.82%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'87.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.53%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6613"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.26%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'30.352.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.31%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'13.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.86%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000007770"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.05%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.08%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'2.135.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'Uniswap"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'5.390"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.29%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.9991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'193.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -4.39%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.184"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.52%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.386"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.50%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'163.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.36%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.944"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.44%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.445"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.292"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.60%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.09129"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.89%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.048"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.82%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.91%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7338"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.91%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.145"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.69%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.706"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.12%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.01793"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.86%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.647"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.06%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.9169"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.48%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.073"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.84%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Quant"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'106.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.67%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'FraxShare"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'5.886"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.29%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.4316"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.35%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.03%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'7.484"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.89%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -5.21%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Aave"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'64.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -11.65%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'NEARProtocol"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1.543"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +6.09%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.018"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.76%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Cronos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.05761"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.50%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Elrond"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'33.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.72%  "
$ws.Range("E51").Style = "Normal"
